$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "-"
$ws.Range("E2").Value = "-"
$ws.Range("C3").Value = "Circuitos Elétricos 2 - MCT-2A"
$ws.Range("D3").Value = "Circuitos Elétricos 2 - MCT-2A"
$ws.Range("C4").Value = "Circuitos Elétricos 2 - ELT-2A"
$ws.Range("C6").Value = "Circuitos Elétricos 2 - ELT-2A"
$ws.Range("E6").Value = "-"
$ws.Range("C7").Value = "-"
$ws.Range("F7").Value = "-"
